$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.747119
$ws.Range("H2").Value = 2.241357
$ws.Range("I2").Value = 0.03096954854571248
$ws.Range("J2").Value = 0.03096954854571248
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 52.91030366666666
$ws.Range("N2").Value = 158.730911
$ws.Range("O2").Value = 0.4161415425564564
$ws.Range("P2").Value = 0.4161415425564564
$ws.Range("Q2").Value = 39.53029316513633
$ws.Range("R2").Value = 355.772638486227
$ws.Range("S2").Value = 0.01288771570408985
$ws.Range("T2").Value = 0.01288771570408985

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.747119
$ws.Range("H3").Value = 2.241357
$ws.Range("I3").Value = 0.03096954854571248
$ws.Range("J3").Value = 0.03096954854571248
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 47.26005833333333
$ws.Range("N3").Value = 141.780175
$ws.Range("O3").Value = 0.3717021489810786
$ws.Range("P3").Value = 0.3717021489810786
$ws.Range("Q3").Value = 35.30888752194166
$ws.Range("R3").Value = 317.779987697475
$ws.Range("S3").Value = 0.01151144774741516
$ws.Range("T3").Value = 0.01151144774741516

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.747119
$ws.Range("H4").Value = 2.241357
$ws.Range("I4").Value = 0.03096954854571248
$ws.Range("J4").Value = 0.03096954854571248
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 26.97460733333333
$ws.Range("N4").Value = 80.923822
$ws.Range("O4").Value = 0.2121563084624651
$ws.Range("P4").Value = 0.2121563084624651
$ws.Range("Q4").Value = 20.15324165627267
$ws.Range("R4").Value = 181.379174906454
$ws.Range("S4").Value = 0.006570385094207465
$ws.Range("T4").Value = 0.006570385094207465

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 19.74619233333334
$ws.Range("H5").Value = 59.23857700000001
$ws.Range("I5").Value = 0.8185184181638298
$ws.Range("J5").Value = 0.8185184181638298
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 52.91030366666666
$ws.Range("N5").Value = 158.730911
$ws.Range("O5").Value = 0.4161415425564564
$ws.Range("P5").Value = 0.4161415425564564
$ws.Range("Q5").Value = 1044.777032617072
$ws.Range("R5").Value = 9402.993293553647
$ws.Range("S5").Value = 0.3406195171455667
$ws.Range("T5").Value = 0.3406195171455667

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 19.74619233333334
$ws.Range("H6").Value = 59.23857700000001
$ws.Range("I6").Value = 0.8185184181638298
$ws.Range("J6").Value = 0.8185184181638298
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 47.26005833333333
$ws.Range("N6").Value = 141.780175
$ws.Range("O6").Value = 0.3717021489810786
$ws.Range("P6").Value = 0.3717021489810786
$ws.Range("Q6").Value = 933.2062015345527
$ws.Range("R6").Value = 8398.855813810975
$ws.Range("S6").Value = 0.3042450550120886
$ws.Range("T6").Value = 0.3042450550120886

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 19.74619233333334
$ws.Range("H7").Value = 59.23857700000001
$ws.Range("I7").Value = 0.8185184181638298
$ws.Range("J7").Value = 0.8185184181638298
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 26.97460733333333
$ws.Range("N7").Value = 80.923822
$ws.Range("O7").Value = 0.2121563084624651
$ws.Range("P7").Value = 0.2121563084624651
$ws.Range("Q7").Value = 532.6457845201439
$ws.Range("R7").Value = 4793.812060681294
$ws.Range("S7").Value = 0.1736538460061745
$ws.Range("T7").Value = 0.1736538460061745

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.630999
$ws.Range("H8").Value = 10.892997
$ws.Range("I8").Value = 0.1505120332904577
$ws.Range("J8").Value = 0.1505120332904577
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 52.91030366666666
$ws.Range("N8").Value = 158.730911
$ws.Range("O8").Value = 0.4161415425564564
$ws.Range("P8").Value = 0.4161415425564564
$ws.Range("Q8").Value = 192.117259703363
$ws.Range("R8").Value = 1729.055337330267
$ws.Range("S8").Value = 0.06263430970679977
$ws.Range("T8").Value = 0.06263430970679978

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.630999
$ws.Range("H9").Value = 10.892997
$ws.Range("I9").Value = 0.1505120332904577
$ws.Range("J9").Value = 0.1505120332904577
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 47.26005833333333
$ws.Range("N9").Value = 141.780175
$ws.Range("O9").Value = 0.3717021489810786
$ws.Range("P9").Value = 0.3717021489810786
$ws.Range("Q9").Value = 171.601224548275
$ws.Range("R9").Value = 1544.411020934475
$ws.Range("S9").Value = 0.05594564622157475
$ws.Range("T9").Value = 0.05594564622157477

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.630999
$ws.Range("H10").Value = 10.892997
$ws.Range("I10").Value = 0.1505120332904577
$ws.Range("J10").Value = 0.1505120332904577
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 26.97460733333333
$ws.Range("N10").Value = 80.923822
$ws.Range("O10").Value = 0.2121563084624651
$ws.Range("P10").Value = 0.2121563084624651
$ws.Range("Q10").Value = 97.944772252726
$ws.Range("R10").Value = 881.502950274534
$ws.Range("S10").Value = 0.03193207736208316
$ws.Range("T10").Value = 0.03193207736208317
